$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.811.17'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.644.90'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  +0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0842'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').Value = '1.870.60'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').Value = '1.651.11'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.17'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.527'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.02%  '
$ws.Range('D17').Value = '26.813.50'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('E18').Value = '  -2.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '213.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.88%  '
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('E21').Value = '  -0.79%  '
$ws.Range('E22').Value = '  +13.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.36'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.20'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.29%  '
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.07'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.69'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.07%  '
$ws.Range('D34').Value = '1.291.19'
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('E37').Value = '  -5.56%  '
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.826'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('E40').Value = '  +0.54%  '
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('D44').Value = '1.796.63'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.53'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0519'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.67'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0979'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.70%  '
